$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.631.82'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.00%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.826.54'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.77%  '

$ws.Range('E4').Value = '  +0.27%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '309.36'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.78%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4680'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +3.68%  '

$ws.Range('E8').Value = '  +0.15%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07130'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.76%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9019'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.11%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07737'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.21%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '19.41'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.25%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.761.96'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.83%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.272'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.17%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.361'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.69%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '87.56'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.18%  '

$ws.Range('E17').Value = '  +0.14%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008546'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.43%  '

$ws.Range('E19').Value = '  +0.25%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '26.673.67'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.11%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.22'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.22%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.022'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.09%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.57'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.37%  '

$ws.Range('E24').Value = '  -3.65%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '153.02'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.32%  '

$ws.Range('E26').Value = '  +0.67%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.980'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.62%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '113.80'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.74%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.867'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.07%  '

$ws.Range('E30').Value = '  +1.65%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.143'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.41%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.815'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.82%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.162'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +5.18%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7372'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.99%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.438'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.04%  '

$ws.Range('E36').Value = '  +1.21%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01933'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.15%  '

$ws.Range('E38').Value = '  +1.44%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.908'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.75%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.883'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.60%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.5051'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.22%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1496'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.11%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.054'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.69%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.008'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.36%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.4664'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.96%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '9.997'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.91%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '97.86'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.17%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.572'
$ws.Range('D48').Style = 'Normal'

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.06047'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.43%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '64.01'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.02%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '35.76'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.54%  '

